$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers (e.g. "353.72")
# need to be kept as TEXT (matching the source data, which stores all Price
# values as strings). Temporarily mark them as Text-formatted before writing
# the value, then restore the default "Normal" style so no stray number
# formatting is left behind on the cell.
$numberLikeDCells = @("D5", "D6", "D9", "D10", "D12", "D13", "D15", "D17", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D43", "D44", "D47", "D49", "D50", "D51")
foreach ($addr in $numberLikeDCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "51.640.10"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").Value = "2.891.28"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "353.72"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").Value = "108.79"
$ws.Range("E6").Value = "  -3.39%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.620"
$ws.Range("E9").Value = "  -2.04%  "

$ws.Range("D10").Value = "38.49"
$ws.Range("E10").Value = "  -4.31%  "

$ws.Range("E11").Value = "  +1.05%  "

$ws.Range("D12").Value = "0.0863"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").Value = "19.32"
$ws.Range("E13").Value = "  -3.37%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.357.32"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "7.66"
$ws.Range("E15").Value = "  -2.04%  "

$ws.Range("D16").Value = "2.898.05"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").Value = "0.968"
$ws.Range("E17").Value = "  -3.97%  "

$ws.Range("D18").Value = "51.650.94"
$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("D19").Value = "3.34"
$ws.Range("E19").Value = "  +0.45%  "

$ws.Range("D20").Value = "7.47"
$ws.Range("E20").Value = "  -2.54%  "

$ws.Range("D21").Value = "13.69"
$ws.Range("E21").Value = "  -3.83%  "

$ws.Range("D22").Value = "0.0₃0970"
$ws.Range("E22").Value = "  -1.24%  "

$ws.Range("D23").Value = "70.09"
$ws.Range("E23").Value = "  -1.44%  "

$ws.Range("D24").Value = "266.45"
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("D26").Value = "0.182"
$ws.Range("E26").Value = "  +8.67%  "

$ws.Range("D27").Value = "26.62"
$ws.Range("E27").Value = "  -0.69%  "

$ws.Range("D28").Value = "7.49"
$ws.Range("E28").Value = "  +14.70%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Value = "0.104"
$ws.Range("E30").Value = "  +8.47%  "

$ws.Range("D31").Value = "10.39"
$ws.Range("E31").Value = "  -2.60%  "

$ws.Range("D32").Value = "36.97"
$ws.Range("E32").Value = "  -2.43%  "

$ws.Range("E33").Value = "  -2.61%  "

$ws.Range("D34").Value = "6.08"
$ws.Range("E34").Value = "  -2.98%  "

$ws.Range("D35").Value = "51.96"
$ws.Range("E35").Value = "  -2.37%  "

$ws.Range("D36").Value = "0.0436"
$ws.Range("E36").Value = "  -3.61%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "3.15"
$ws.Range("E38").Value = "  -5.64%  "

$ws.Range("D39").Value = "18.02"
$ws.Range("E39").Value = "  -4.36%  "

$ws.Range("D40").Value = "1.98"
$ws.Range("E40").Value = "  -4.73%  "

$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  -7.13%  "

$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("D43").Value = "22.77"
$ws.Range("E43").Value = "  -4.09%  "

$ws.Range("D44").Value = "118.67"
$ws.Range("E44").Value = "  -2.68%  "

$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("E46").Value = "  -5.58%  "

$ws.Range("D47").Value = "3.40"
$ws.Range("E47").Value = "  -4.53%  "

$ws.Range("D48").Value = "2.115.44"
$ws.Range("E48").Value = "  -3.84%  "

$ws.Range("D49").Value = "0.247"
$ws.Range("E49").Value = "  -7.01%  "

$ws.Range("D50").Value = "0.0335"
$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").Value = "9.01"
$ws.Range("E51").Value = "  -0.91%  "

foreach ($addr in $numberLikeDCells) {
    $ws.Range($addr).Style = "Normal"
}

